$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 ---
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 34.080206910000001
$ws.Range("D2").Value = 12.675873040000001
$ws.Range("E2").Value = 19.31180367
$ws.Range("F2").Value = 2.2248757
$ws.Range("G2").Value = 22.80563558
$ws.Range("H2").Value = 2.42225743
$ws.Range("I2").Value = 5.21670989
$ws.Range("J2").Value = 0.48746172

# --- Update row 3 ---
$ws.Range("B3").Value = 2
$ws.Range("C3").Value = 25.99566077
$ws.Range("D3").Value = 9.0066983199999999
$ws.Range("E3").Value = 16.74555297
$ws.Range("F3").Value = 4.38207694
$ws.Range("G3").Value = 15.31025662
$ws.Range("H3").Value = 2.50698321
$ws.Range("I3").Value = 5.01860437
$ws.Range("J3").Value = 0.33179396

# --- Row 4 ---
$ws.Range("A4").Value = "cross2"
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = 19.60745704
$ws.Range("D4").Value = 8.41630748
$ws.Range("E4").Value = 14.55165423
$ws.Range("F4").Value = 3.29758928
$ws.Range("G4").Value = 9.7910599
$ws.Range("H4").Value = 2.38047576
$ws.Range("I4").Value = 5.07091151
$ws.Range("J4").Value = 0.48098059

# --- Row 5 ---
$ws.Range("A5").Value = "cross2"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = 16.495112
$ws.Range("D5").Value = 4.13411092
$ws.Range("E5").Value = 13.07396301
$ws.Range("F5").Value = 1.82039157
$ws.Range("G5").Value = 7.65384726
$ws.Range("H5").Value = 2.01195022
$ws.Range("I5").Value = 5.62253413
$ws.Range("J5").Value = 0.80904396

# --- Row 6 ---
$ws.Range("A6").Value = "cross2"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = 12.56025587
$ws.Range("D6").Value = 2.67137785
$ws.Range("E6").Value = 12.26705992
$ws.Range("F6").Value = 2.28747686
$ws.Range("G6").Value = 6.35174595
$ws.Range("H6").Value = 1.71292113
$ws.Range("I6").Value = 5.44861621
$ws.Range("J6").Value = 0.69452316

# --- Row 7 ---
$ws.Range("A7").Value = "cross2"
$ws.Range("B7").Value = 6
$ws.Range("C7").Value = 10.63080093
$ws.Range("D7").Value = 2.18667672
$ws.Range("E7").Value = 10.62519117
$ws.Range("F7").Value = 2.93331336
$ws.Range("G7").Value = 5.80242588
$ws.Range("H7").Value = 0.86696993
$ws.Range("I7").Value = 5.56298189
$ws.Range("J7").Value = 0.64764543

# --- Row 8 ---
$ws.Range("A8").Value = "cross2"
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 8.76711844
$ws.Range("D8").Value = 0.98441886
$ws.Range("E8").Value = 8.22543453
$ws.Range("F8").Value = 1.32651864
$ws.Range("G8").Value = 5.1452006
$ws.Range("H8").Value = 0.67946208
$ws.Range("I8").Value = 5.11745961
$ws.Range("J8").Value = 0.83527476

# --- Row 9 ---
$ws.Range("A9").Value = "cross2"
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = 8.22562096
$ws.Range("D9").Value = 1.79769348
$ws.Range("E9").Value = 6.97149843
$ws.Range("F9").Value = 2.41071723
$ws.Range("G9").Value = 4.81346085
$ws.Range("H9").Value = 0.83370028
$ws.Range("I9").Value = 4.66684735
$ws.Range("J9").Value = 0.91118546

# --- Row 10 ---
$ws.Range("A10").Value = "cross2"
$ws.Range("B10").Value = 20
$ws.Range("C10").Value = 6.67888612
$ws.Range("D10").Value = 0.70903374
$ws.Range("E10").Value = 4.75171762
$ws.Range("F10").Value = 1.66211799
$ws.Range("G10").Value = 3.87561651
$ws.Range("H10").Value = 0.84365344
$ws.Range("I10").Value = 3.84293909
$ws.Range("J10").Value = 0.92278434

# --- Selection ---
[void]$ws.Range("C15").Select()

Write-Output "done"
